$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace technician name in H3 (was "Dakota Myers") with "Franz Ferdinand"
$ws.Range("H3").Value = "Franz Ferdinand"

# Reflect the selected cell as recorded at save time
$ws.Range("H3").Select()
